$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '90.148.65'
$ws.Range('E2').Value = '  +2.48%  '
$ws.Range('D3').Value = '3.187.93'
$ws.Range('E3').Value = '  -2.01%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = '''214.67'
$ws.Range('E5').Value = '  +0.64%  '
$ws.Range('D6').Value = '''619.64'
$ws.Range('E6').Value = '  -1.51%  '
$ws.Range('D7').Value = '''0.398'
$ws.Range('E7').Value = '  +2.59%  '
$ws.Range('D8').Value = '''0.692'
$ws.Range('E8').Value = '  -2.77%  '
$ws.Range('E9').Value = '  +0.02%  '
$ws.Range('D10').Value = '3.182.79'
$ws.Range('E10').Value = '  -1.99%  '
$ws.Range('D11').Value = '''0.579'
$ws.Range('E11').Value = '  +0.02%  '
$ws.Range('E12').Value = '  -6.13%  '
$ws.Range('D13').Value = '''0.0000257'
$ws.Range('E13').Value = '  -4.34%  '
$ws.Range('D14').Value = '89.982.38'
$ws.Range('E14').Value = '  +2.62%  '
$ws.Range('D15').Value = '3.774.45'
$ws.Range('E15').Value = '  -2.07%  '
$ws.Range('D16').Value = '''33.03'
$ws.Range('E16').Value = '  -3.89%  '
$ws.Range('D17').Value = '''5.27'
$ws.Range('E17').Value = '  -4.06%  '
$ws.Range('D18').Value = '3.180.87'
$ws.Range('E18').Value = '  -2.61%  '
$ws.Range('E19').Value = '  +0.56%  '
$ws.Range('D20').Value = '''0.0000204'
$ws.Range('E20').Value = '  +42.69%  '
$ws.Range('D21').Value = '''13.43'
$ws.Range('E21').Value = '  -4.68%  '
$ws.Range('D22').Value = '''439.37'
$ws.Range('E22').Value = '  +0.58%  '
$ws.Range('D23').Value = '''8.63'
$ws.Range('E23').Value = '  -4.23%  '
$ws.Range('D24').Value = '''5.08'
$ws.Range('E24').Value = '  -4.95%  '
$ws.Range('D25').Value = '''5.16'
$ws.Range('E25').Value = '  -4.02%  '
$ws.Range('D26').Value = '''11.66'
$ws.Range('E26').Value = '  -6.61%  '
$ws.Range('D27').Value = '3.348.11'
$ws.Range('E27').Value = '  -2.39%  '
$ws.Range('D28').Value = '''75.43'
$ws.Range('E28').Value = '  -2.37%  '
$ws.Range('E29').Value = '  +0.36%  '
$ws.Range('E30').Value = '  -5.90%  '
$ws.Range('E31').Value = '  -0.05%  '
$ws.Range('D32').Value = '''4.18'
$ws.Range('E32').Value = '  +27.05%  '
$ws.Range('D33').Value = '''8.45'
$ws.Range('E33').Value = '  -5.03%  '
$ws.Range('D34').Value = '''537.19'
$ws.Range('E34').Value = '  -5.87%  '
$ws.Range('E35').Value = '  -3.08%  '
$ws.Range('E36').Value = '  -5.45%  '
$ws.Range('D37').Value = '''1.28'
$ws.Range('E37').Value = '  -8.33%  '
$ws.Range('D38').Value = '''22.08'
$ws.Range('E38').Value = '  -4.03%  '
$ws.Range('D39').Value = '''22.32'
$ws.Range('E39').Value = '  +2.34%  '
$ws.Range('E40').Value = '  -8.82%  '
$ws.Range('D41').Value = '''0.998'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('B42').Value = 'Stacks'
$ws.Range('C42').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D42').Value = '''1.95'
$ws.Range('E42').Value = '  -4.70%  '
$ws.Range('B43').Value = 'USDe'
$ws.Range('C43').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D43').Value = '''1.00'
$ws.Range('E43').Value = '  +0.13%  '
$ws.Range('D44').Value = '''0.375'
$ws.Range('E44').Value = '  -7.19%  '
$ws.Range('D45').Value = '''150.07'
$ws.Range('E45').Value = '  -1.04%  '
$ws.Range('D46').Value = '''43.81'
$ws.Range('E46').Value = '  -3.20%  '
$ws.Range('D47').Value = '''172.55'
$ws.Range('E47').Value = '  -4.44%  '
$ws.Range('E48').Value = '  -8.98%  '
$ws.Range('E49').Value = '  -8.26%  '
$ws.Range('D50').Value = '''4.06'
$ws.Range('E50').Value = '  -4.75%  '
$ws.Range('D51').Value = '''0.611'
$ws.Range('E51').Value = '  -3.84%  '
